# Generate Report for Handoff
# Adds two newly-handed-off files (a65632bd-... and c566f324-...) to the
# localization-status workbook: one new row on each of the three sheets
# (Overview / zh-cn / de-de), each carrying its own hyperlinks.

$wb = $excel.ActiveWorkbook

$uuid1 = "a65632bd-c226-4dfb-8e2d-40917d54d4b9"
$uuid2 = "c566f324-9411-43c6-bb0c-4790b242ec09"

$xlf1ZhCn = "a65632bd-c226-4dfb-8e2d-40917d54d4b9.27a9f3d3f64f12675bddeccbeabf7328e997b20c.zh-cn.xlf"
$xlf2ZhCn = "c566f324-9411-43c6-bb0c-4790b242ec09.5327ae33e62c4a7333911e84799207e2cf8f6baa.zh-cn.xlf"
$xlf1DeDe = "a65632bd-c226-4dfb-8e2d-40917d54d4b9.27a9f3d3f64f12675bddeccbeabf7328e997b20c.de-de.xlf"
$xlf2DeDe = "c566f324-9411-43c6-bb0c-4790b242ec09.5327ae33e62c4a7333911e84799207e2cf8f6baa.de-de.xlf"

$statusReady = "Ready for handoff"
$dotMd = ".md"
$include = "Include"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1): File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(6, 1), "https://github.com/OpenLocalizationTest/oltest/blob/d1b74bdc5d3a66400dc64345f463d4e92aa689c8/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$wsOverview.Cells.Item(6, 2).Value = $statusReady
$wsOverview.Cells.Item(6, 3).Value = $statusReady
$wsOverview.Cells.Item(6, 4).Value = "2016-36-11 10:36:14"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(7, 1), "https://github.com/OpenLocalizationTest/oltest/blob/4fa5856656fc54d9066abd67807646764148a607/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$wsOverview.Cells.Item(7, 2).Value = $statusReady
$wsOverview.Cells.Item(7, 3).Value = $statusReady
$wsOverview.Cells.Item(7, 4).Value = "2016-36-11 10:36:14"

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(6, 1), "https://github.com/OpenLocalizationTest/oltest/blob/d1b74bdc5d3a66400dc64345f463d4e92aa689c8/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(6, 2), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b5f8821a6ccf6a2825cf504556636af8f22de16d/e2e/$uuid1.md", "", "", $dotMd) | Out-Null
$wsZhCn.Cells.Item(6, 3).Value = $statusReady
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(6, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/338e1c715277d3557ec73993f3f735112ca23adb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf1ZhCn", "", "", $xlf1ZhCn) | Out-Null
$wsZhCn.Cells.Item(6, 5).Value = "2016-03-11 10:36:11"
$wsZhCn.Cells.Item(6, 8).Value = $epoch
$wsZhCn.Cells.Item(6, 9).Value = $include

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(7, 1), "https://github.com/OpenLocalizationTest/oltest/blob/4fa5856656fc54d9066abd67807646764148a607/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(7, 2), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1db30bc0ae0069abd7a4fd8ac052e842dc9ec9f2/e2e/$uuid2.md", "", "", $dotMd) | Out-Null
$wsZhCn.Cells.Item(7, 3).Value = $statusReady
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(7, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/58737261fbc31b879ab517ce3ba53e89efb08041/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf2ZhCn", "", "", $xlf2ZhCn) | Out-Null
$wsZhCn.Cells.Item(7, 5).Value = "2016-03-11 10:36:11"
$wsZhCn.Cells.Item(7, 8).Value = $epoch
$wsZhCn.Cells.Item(7, 9).Value = $include

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(6, 1), "https://github.com/OpenLocalizationTest/oltest/blob/d1b74bdc5d3a66400dc64345f463d4e92aa689c8/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(6, 2), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3d906a94a6ba25df19db8fe27d5aa2f0c67269d4/e2e/$uuid1.md", "", "", $dotMd) | Out-Null
$wsDeDe.Cells.Item(6, 3).Value = $statusReady
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(6, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee871baa0e0f904e93dc4dde6e37c49c5fafde55/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf1DeDe", "", "", $xlf1DeDe) | Out-Null
$wsDeDe.Cells.Item(6, 5).Value = "2016-03-11 10:36:14"
$wsDeDe.Cells.Item(6, 8).Value = $epoch
$wsDeDe.Cells.Item(6, 9).Value = $include

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(7, 1), "https://github.com/OpenLocalizationTest/oltest/blob/4fa5856656fc54d9066abd67807646764148a607/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(7, 2), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/443620fc16ee2eb0dfc4b33ee722e7a6fdc85a10/e2e/$uuid2.md", "", "", $dotMd) | Out-Null
$wsDeDe.Cells.Item(7, 3).Value = $statusReady
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(7, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2810a132bd0aa061422465b49c4cc3d0093f90c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf2DeDe", "", "", $xlf2DeDe) | Out-Null
$wsDeDe.Cells.Item(7, 5).Value = "2016-03-11 10:36:14"
$wsDeDe.Cells.Item(7, 8).Value = $epoch
$wsDeDe.Cells.Item(7, 9).Value = $include

Write-Host "Added handoff rows for $uuid1 and $uuid2 across Overview/zh-cn/de-de sheets."
